# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (interested-count) figures in column F for two rows
# on both the "展览" sheet and the "全部类型" sheet (which mirrors it).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1858
$ws1.Range("F6").Value = 1114

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1858
$ws4.Range("F6").Value = 1115
